# Generate Report for Handoff
# Adds two new handed-off files (8d0cffe5-... and c61d2558-...) to the
# localization-status workbook: one new row on "Overview" and one new row
# on each of the "zh-cn" / "de-de" detail sheets, including their hyperlinks.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A6").Value = "8d0cffe5-57db-4a96-890b-cf2d0d0dd30d.md"
$wsOverview.Range("B6").Value = "Ready for handoff"
$wsOverview.Range("C6").Value = "Ready for handoff"
$wsOverview.Range("D6").Value = "2016-32-18 22:32:46"

$wsOverview.Range("A7").Value = "c61d2558-f073-4f61-9836-6cf826f42d0b.md"
$wsOverview.Range("B7").Value = "Ready for handoff"
$wsOverview.Range("C7").Value = "Ready for handoff"
$wsOverview.Range("D7").Value = "2016-32-18 22:32:46"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/8d0cffe5-57db-4a96-890b-cf2d0d0dd30d/e2e/8d0cffe5-57db-4a96-890b-cf2d0d0dd30d.md", "", "", "8d0cffe5-57db-4a96-890b-cf2d0d0dd30d.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/c61d2558-f073-4f61-9836-6cf826f42d0b/e2e/c61d2558-f073-4f61-9836-6cf826f42d0b.md", "", "", "c61d2558-f073-4f61-9836-6cf826f42d0b.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A6").Value = "8d0cffe5-57db-4a96-890b-cf2d0d0dd30d.md"
$wsZhCn.Range("B6").Value = ".md"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("D6").Value = "8d0cffe5-57db-4a96-890b-cf2d0d0dd30d.9af80e1b8762b04e1624d994f4c8047f472315c7.zh-cn.xlf"
$wsZhCn.Range("E6").Value = "2016-03-18 22:32:43"
$wsZhCn.Range("E6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H6").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I6").Value = "Include"

$wsZhCn.Range("A7").Value = "c61d2558-f073-4f61-9836-6cf826f42d0b.md"
$wsZhCn.Range("B7").Value = ".md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("D7").Value = "c61d2558-f073-4f61-9836-6cf826f42d0b.fc21d56db6f7cf8d1a258ba89ae1ae775c17f339.zh-cn.xlf"
$wsZhCn.Range("E7").Value = "2016-03-18 22:32:43"
$wsZhCn.Range("E7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H7").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I7").Value = "Include"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/8d0cffe5-57db-4a96-890b-cf2d0d0dd30d/e2e/8d0cffe5-57db-4a96-890b-cf2d0d0dd30d.md", "", "", "8d0cffe5-57db-4a96-890b-cf2d0d0dd30d.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B6"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/8d0cffe5-57db-4a96-890b-cf2d0d0dd30d/e2e/8d0cffe5-57db-4a96-890b-cf2d0d0dd30d.md", "", "", ".md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8d0cffe5-57db-4a96-890b-cf2d0d0dd30d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8d0cffe5-57db-4a96-890b-cf2d0d0dd30d.9af80e1b8762b04e1624d994f4c8047f472315c7.zh-cn.xlf", "", "", "8d0cffe5-57db-4a96-890b-cf2d0d0dd30d.9af80e1b8762b04e1624d994f4c8047f472315c7.zh-cn.xlf") | Out-Null

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/c61d2558-f073-4f61-9836-6cf826f42d0b/e2e/c61d2558-f073-4f61-9836-6cf826f42d0b.md", "", "", "c61d2558-f073-4f61-9836-6cf826f42d0b.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B7"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c61d2558-f073-4f61-9836-6cf826f42d0b/e2e/c61d2558-f073-4f61-9836-6cf826f42d0b.md", "", "", ".md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c61d2558-f073-4f61-9836-6cf826f42d0b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c61d2558-f073-4f61-9836-6cf826f42d0b.fc21d56db6f7cf8d1a258ba89ae1ae775c17f339.zh-cn.xlf", "", "", "c61d2558-f073-4f61-9836-6cf826f42d0b.fc21d56db6f7cf8d1a258ba89ae1ae775c17f339.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A6").Value = "8d0cffe5-57db-4a96-890b-cf2d0d0dd30d.md"
$wsDeDe.Range("B6").Value = ".md"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("D6").Value = "8d0cffe5-57db-4a96-890b-cf2d0d0dd30d.9af80e1b8762b04e1624d994f4c8047f472315c7.de-de.xlf"
$wsDeDe.Range("E6").Value = "2016-03-18 22:32:46"
$wsDeDe.Range("E6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H6").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I6").Value = "Include"

$wsDeDe.Range("A7").Value = "c61d2558-f073-4f61-9836-6cf826f42d0b.md"
$wsDeDe.Range("B7").Value = ".md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("D7").Value = "c61d2558-f073-4f61-9836-6cf826f42d0b.fc21d56db6f7cf8d1a258ba89ae1ae775c17f339.de-de.xlf"
$wsDeDe.Range("E7").Value = "2016-03-18 22:32:46"
$wsDeDe.Range("E7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H7").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I7").Value = "Include"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/8d0cffe5-57db-4a96-890b-cf2d0d0dd30d/e2e/8d0cffe5-57db-4a96-890b-cf2d0d0dd30d.md", "", "", "8d0cffe5-57db-4a96-890b-cf2d0d0dd30d.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B6"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8d0cffe5-57db-4a96-890b-cf2d0d0dd30d/e2e/8d0cffe5-57db-4a96-890b-cf2d0d0dd30d.md", "", "", ".md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8d0cffe5-57db-4a96-890b-cf2d0d0dd30d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8d0cffe5-57db-4a96-890b-cf2d0d0dd30d.9af80e1b8762b04e1624d994f4c8047f472315c7.de-de.xlf", "", "", "8d0cffe5-57db-4a96-890b-cf2d0d0dd30d.9af80e1b8762b04e1624d994f4c8047f472315c7.de-de.xlf") | Out-Null

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/c61d2558-f073-4f61-9836-6cf826f42d0b/e2e/c61d2558-f073-4f61-9836-6cf826f42d0b.md", "", "", "c61d2558-f073-4f61-9836-6cf826f42d0b.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B7"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/c61d2558-f073-4f61-9836-6cf826f42d0b/e2e/c61d2558-f073-4f61-9836-6cf826f42d0b.md", "", "", ".md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c61d2558-f073-4f61-9836-6cf826f42d0b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c61d2558-f073-4f61-9836-6cf826f42d0b.fc21d56db6f7cf8d1a258ba89ae1ae775c17f339.de-de.xlf", "", "", "c61d2558-f073-4f61-9836-6cf826f42d0b.fc21d56db6f7cf8d1a258ba89ae1ae775c17f339.de-de.xlf") | Out-Null

Write-Host "Handback report rows added."
